# Generate Report for Handoff
# Rewrites the localization-status report for the newly-generated source
# file (032ceea2-... -> 34f6c38c-...) and its freshly re-handed-off
# xliff files, resetting the (now stale) handback info.

$wb = $excel.ActiveWorkbook

$oldGuid = "032ceea2-d288-4b97-9cb2-4850214ad6e1"
$newGuid = "34f6c38c-717b-4187-8b13-f872ad69a5ab"
$newHash = "6493fc83ded4c0c70a3cd6fbb908151c3643df43"

$newMdName   = "$newGuid.md"
$newMdPath   = "e2e\$newGuid.md"
$newHoDate   = "2016-08-21 23:07:01"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh       = $wb.Worksheets.Item("zh-cn")
$wsDe       = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdPath
$wsOverview.Range("G2").Value = $newHoDate

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = $newMdPath
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-21 23:06:56"
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMdName
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = $newHoDate
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMdName
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
